$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.00942850112915
$ws.Range("B1").Value = 1.348303079605103
$ws.Range("C1").Value = 2.179309844970703
$ws.Range("D1").Value = 4.447084426879883
$ws.Range("E1").Value = 2.048956155776978
